{"js": "// Replace the unfilled \"<invullen>\" student-number placeholder with the\n// actual student number (851941098) \u2014 see commit \"adding student number\n// randy\".\n//\n// The placeholder paragraph originally reads \"<invullen>\": an italic run\n// \"<\", an italic run \"invullen\", a (collapsed) \"_GoBack\" bookmark, then an\n// italic run \">\". We find that exact paragraph, swap its text for the\n// student number (which also collapses the runs into one, inheriting the\n// paragraph's italic/nl-NL formatting), and then restore the \"_GoBack\"\n// bookmark at the start of the paragraph so it keeps marking this edit\n// location, just like in the original file.\n\nconst placeholder = \"<invullen>\";\nconst studentNumber = \"851941098\";\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfunction findParagraphByText(items, text) {\n  for (let i = 0; i < items.length; i++) {\n    if (items[i].text === text) return items[i];\n  }\n  return null;\n}\n\nlet target = findParagraphByText(paragraphs.items, placeholder);\n\nif (target) {\n  // Replace the whole paragraph's text in one go (scoped to this single\n  // paragraph, so the similar \"<...>\" placeholder later in the report is\n  // left untouched).\n  target.getRange().insertText(studentNumber, \"Replace\");\n  await context.sync();\n} else {\n  // Fallback: locate the exact placeholder text anywhere in the body.\n  const matches = body.search(placeholder, { matchCase: true });\n  matches.load(\"items\");\n  await context.sync();\n  if (matches.items.length > 0) {\n    matches.items[0].insertText(studentNumber, \"Replace\");\n    await context.sync();\n  }\n}\n\n// Re-seat the \"_GoBack\" bookmark at the start of the (now renumbered)\n// paragraph, matching its original collapsed position in that paragraph.\nconst paragraphs2 = context.document.body.paragraphs;\nparagraphs2.load(\"items/text\");\nawait context.sync();\nconst target2 = findParagraphByText(paragraphs2.items, studentNumber);\nif (target2) {\n  target2.getRange(\"Start\").insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# Replace the unfilled \"<invullen>\" student-number placeholder with the\n# actual student number (851941098) - see commit \"adding student number\n# randy\".\n#\n# The placeholder paragraph originally reads \"<invullen>\": an italic run\n# \"<\", an italic run \"invullen\", a (collapsed) \"_GoBack\" bookmark, then an\n# italic run \">\". We find that exact paragraph, swap its text for the\n# student number (which also collapses the runs into one, inheriting the\n# paragraph's italic/nl-NL formatting), and then restore the \"_GoBack\"\n# bookmark at the start of the paragraph so it keeps marking this edit\n# location, just like in the original file.\n\n$d = $word.ActiveDocument\n\n$placeholder = \"<invullen>\"\n$studentNumber = \"851941098\"\n\n$target = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $para = $d.Paragraphs.Item($i)\n    if ($para.Range.Text -eq $placeholder) {\n        $target = $para\n        break\n    }\n}\n\nif ($target -ne $null) {\n    $r = $target.Range\n    # Drop the trailing paragraph mark from the range so only the visible\n    # text is replaced.\n    $r.MoveEnd(1, -1) | Out-Null\n    $r.Text = $studentNumber\n} else {\n    # Fallback: scoped Find/Replace for the exact placeholder text only.\n    $find = $d.Content.Find\n    $find.Text = $placeholder\n    $find.Replacement.Text = $studentNumber\n    $find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $studentNumber, 2) | Out-Null\n}\n\n# Re-seat the \"_GoBack\" bookmark at the start of the (now renumbered)\n# paragraph, matching its original collapsed position in that paragraph.\n$target2 = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $para = $d.Paragraphs.Item($i)\n    if ($para.Range.Text -eq $studentNumber) {\n        $target2 = $para\n        break\n    }\n}\n\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\nif ($target2 -ne $null) {\n    $startPos = $target2.Range.Start\n    $bmRange = $d.Range($startPos, $startPos)\n    $d.Bookmarks.Add(\"_GoBack\", $bmRange) | Out-Null\n}\n"}
